# Weekly fruit/vegetable price update: insert a new record as row 48
# (Agrícola del Norte S.A. de Arica - Uva), shifting all subsequent rows
# down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 48; existing rows 48..114 shift to 49..115
$ws.Rows.Item(48).Insert()

# Populate the newly inserted row 48 with the new price record.
$ws.Range("A48").Value = 1
$ws.Range("B48").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C48").Value = "Arica y Parinacota"
$ws.Range("D48").Value = 44914
$ws.Range("E48").Value = 15
$ws.Range("F48").Value = "Fruta"
$ws.Range("G48").Value = 100109
$ws.Range("H48").Value = "Uva"
$ws.Range("I48").Value = 100109001
$ws.Range("J48").Value = "Uva"
$ws.Range("K48").Value = "Superior Seedless"
$ws.Range("L48").Value = "Segunda"
$ws.Range("M48").Value = 450
$ws.Range("N48").Value = 14000
$ws.Range("O48").Value = 15000
$ws.Range("P48").Value = 14556
$ws.Range("Q48").Value = "$/bandeja 10 kilos"
$ws.Range("R48").Value = "Región de Coquimbo"
$ws.Range("S48").Value = 1456
$ws.Range("T48").Value = 10
